# Auto-update data + news
# Update the "Initial Jobless Claims" (ICSA_thou) row (row 9) with refreshed
# values: the reported value, its 10-period average, and the absolute/percent
# change versus the prior reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value = 198000
$ws.Range("G9").Value = 364318.0076628352
$ws.Range("H9").Value = -7000
$ws.Range("I9").Value = -0.03414634146341464
